$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.781.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "'3.761.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'424.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.90%  "

$ws.Range("D6").Value = "'137.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.89%  "

$ws.Range("D7").Value = "'0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'0.725"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.06%  "

$ws.Range("E10").Value = "  -9.65%  "

$ws.Range("D11").Value = "'0.0000297"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.70%  "

$ws.Range("D12").Value = "'42.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.08%  "

$ws.Range("D13").Value = "'10.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.55%  "

$ws.Range("D14").Value = "'4.370.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Value = "'14.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.16%  "

$ws.Range("D16").Value = "'0.137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").Value = "'3.750.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "'19.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.73%  "

$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").Value = "'65.929.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").Value = "'399.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.80%  "

$ws.Range("D22").Value = "'14.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.36%  "

$ws.Range("D23").Value = "'3.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.06%  "

$ws.Range("D24").Value = "'83.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.19%  "

$ws.Range("D25").Value = "'36.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").Value = "'9.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +39.16%  "

$ws.Range("D27").Value = "'3.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.46%  "

$ws.Range("E28").Value = "  +4.68%  "

$ws.Range("E29").Value = "  -3.61%  "

$ws.Range("D30").Value = "'13.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.42%  "

$ws.Range("D31").Value = "'699.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("D32").Value = "'0.131"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.77%  "

$ws.Range("E33").Value = "  +1.72%  "

$ws.Range("D34").Value = "'40.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.86%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'5.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +32.94%  "

$ws.Range("D37").Value = "'56.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.48%  "

$ws.Range("D38").Value = "'0.147"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.56%  "

$ws.Range("D39").Value = "'0.0465"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.68%  "

$ws.Range("E40").Value = "  +42.34%  "

$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "'0.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.06%  "

$ws.Range("D44").Value = "'0.0₃0657"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.37%  "

$ws.Range("D45").Value = "'3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "

$ws.Range("D46").Value = "'3.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.15%  "

$ws.Range("D47").Value = "'0.316"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.98%  "

$ws.Range("E48").Value = "  +4.74%  "

$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").Value = "'140.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.11%  "

$ws.Range("D51").Value = "'2.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.10%  "
